$d = $word.ActiveDocument

# --- Change 1: "Visual Studio – " -> "Visual Studio" + " Code" + " – " (3 runs) ---
$find = $d.Content
$found = $find.Find.Execute("Visual Studio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertAt = $d.Range($find.End, $find.End)
    $insertAt.InsertAfter(" Code")
    $newRng = $d.Range($find.End, $find.End + 5)
    # Force a genuine run split by toggling a character property, then
    # reverting it so the final formatting matches the surrounding text.
    $newRng.Font.Bold = $true
    $newRng.Font.Bold = $false
}

# --- Change 2: remove the PowerPoint / Excel / trailing empty-paragraph rows ---
# (the "Microsoft PowerPoint – Presentation" and "Microsoft Excel – QA
# Documentation" bullet paragraphs, plus the blank paragraph right after
# them, are deleted outright; the paragraphs before and after survive.)
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Microsoft PowerPoint*") {
        $startIndex = $i
        break
    }
}

if ($startIndex -ge 1) {
    $endIndex = [Math]::Min($startIndex + 2, $d.Paragraphs.Count)
    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($endIndex)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
